$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value2 = 'filename'
$ws.Range("B1").Value2 = 'name'
$ws.Range("C1").Value2 = 'surname'
$ws.Range("D1").Value2 = 'university'
$ws.Range("E1").Value2 = 'bank account'
$ws.Range("F1").Value2 = 'currency'
$ws.Range("G1").Value2 = 'amount'
$ws.Range("H1").Value2 = 'case number'

$ws.Range("A2").Value2 = 'Ella_Allen_University of Opole.xlsx'
$ws.Range("B2").Value2 = 'Ella'
$ws.Range("C2").Value2 = 'Allen'
$ws.Range("D2").Value2 = 'University of Opole'
$ws.Range("E2").Value2 = 'HU12345678901234567890123456'
$ws.Range("F2").Value2 = 'HUF'
$ws.Range("G2").Value2 = 3456789.01
$ws.Range("H2").Value2 = 'XLxFodWStBq9vqp'

$ws.Range("A3").Value2 = 'Fiona_Garcia_AGH University of Science and Technology.xlsx'
$ws.Range("B3").Value2 = 'Fiona'
$ws.Range("C3").Value2 = 'Garcia'
$ws.Range("D3").Value2 = 'AGH University of Science and Technology'
$ws.Range("E3").Value2 = 'AU90123456789012345678901234'
$ws.Range("F3").Value2 = 'AUD'
$ws.Range("G3").Value2 = 2109.87
$ws.Range("H3").Value2 = 'TSSqJTgQscU3xRk'

$ws.Range("A4").Value2 = 'Jack_Adams_University of Innsbruck.xlsx'
$ws.Range("B4").Value2 = 'Jack'
$ws.Range("C4").Value2 = 'Adams'
$ws.Range("D4").Value2 = 'University of Innsbruck'
$ws.Range("E4").Value2 = 'CH56789012345678901234567890'
$ws.Range("F4").Value2 = 'CHF'
$ws.Range("G4").Value2 = 8765.43
$ws.Range("H4").Value2 = '20PHrLFFPkjoFh3'

$ws.Range("A5").Value2 = 'Kevin_Hernandez_Silesian University of Technology.xlsx'
$ws.Range("B5").Value2 = 'Kevin'
$ws.Range("C5").Value2 = 'Hernandez'
$ws.Range("D5").Value2 = 'Silesian University of Technology'
$ws.Range("E5").Value2 = 'DK56789012345678901234567890'
$ws.Range("F5").Value2 = 'DKK'
$ws.Range("G5").Value2 = 54321.09
$ws.Range("H5").Value2 = 'IntCHCcO5hVBx7F'

$ws.Range("A6").Value2 = 'Olivia_Roberts_University of Bern.xlsx'
$ws.Range("B6").Value2 = 'Olivia'
$ws.Range("C6").Value2 = 'Roberts'
$ws.Range("D6").Value2 = 'University of Bern'
$ws.Range("E6").Value2 = 'JP12345678901234567890'
$ws.Range("F6").Value2 = 'JPY'
$ws.Range("G6").Value2 = 876543.21
$ws.Range("H6").Value2 = 'PsvOIOwtWHkuhak'

$ws.Range("A7").Value2 = 'Patricia_Thomas_University of Silesia.xlsx'
$ws.Range("B7").Value2 = 'Patricia'
$ws.Range("C7").Value2 = 'Thomas'
$ws.Range("D7").Value2 = 'University of Silesia'
$ws.Range("E7").Value2 = 'PT12345678901234567890123456'
$ws.Range("F7").Value2 = 'EUR'
$ws.Range("G7").Value2 = 7654.32
$ws.Range("H7").Value2 = 'azdzUIGGaiGBDMW'

$ws.Range("A8").Value2 = 'Ursula_Lee_Opole University.xlsx'
$ws.Range("B8").Value2 = 'Ursula'
$ws.Range("C8").Value2 = 'Lee'
$ws.Range("D8").Value2 = 'Opole University'
$ws.Range("E8").Value2 = 'LU90123456789012345678901234'
$ws.Range("F8").Value2 = 'EUR'
$ws.Range("G8").Value2 = 9876.54
$ws.Range("H8").Value2 = 'RUk3wecSSbPhWC0'

$ws.Range("A9").Value2 = 'Zane_Clark_University of Economics in Katowice.xlsx'
$ws.Range("B9").Value2 = 'Zane'
$ws.Range("C9").Value2 = 'Clark'
$ws.Range("D9").Value2 = 'University of Economics in Katowice'
$ws.Range("E9").Value2 = 'EE56789012345678901234567890'
$ws.Range("F9").Value2 = 'EUR'
$ws.Range("G9").Value2 = 6543.21
$ws.Range("H9").Value2 = 'Ki1E2GL7fuVqd1n'

$ws.Rows.Item(10).ClearContents()

$ws.Columns.Item(2).ColumnWidth = 9.17
$ws.Columns.Item(6).ColumnWidth = 9.17
